$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.740.44'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '1.547.20'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'206.18"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = "'21.48"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.23%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = "'0.0582"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').Value = "'0.0853"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('D12').Value = '1.767.01'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '1.542.82'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = "'3.68"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = "'0.512"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '26.733.48'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('D17').Value = "'61.29"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').Value = "'213.12"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = '0.0₃0689'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = "'7.24"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').Value = "'8.96"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('D25').Value = "'152.52"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = "'6.50"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.63%  '
$ws.Range('D27').Value = "'14.88"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  -0.85%  '
$ws.Range('D30').Value = "'0.0461"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('D32').Value = "'3.20"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').Value = '1.344.35'
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('D34').Value = "'2.92"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('D37').Value = "'0.937"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = "'0.524"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('D40').Value = "'5.82"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.45%  '
$ws.Range('D41').Value = "'0.801"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('E42').Value = '  -1.21%  '
$ws.Range('D43').Value = "'2.20"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = "'62.75"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('E45').Value = '  -4.08%  '
$ws.Range('D46').Value = '1.680.95'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('E47').Value = '  -3.95%  '
$ws.Range('D48').Value = "'86.02"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('D50').Value = '0.0₇0974'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('D51').Value = "'0.0953"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.34%  '
